$wb = $excel.ActiveWorkbook

# ---- Sheet "main" ----
$ws1 = $wb.Worksheets.Item("main")

# Header row
$ws1.Range("B1").Value = "Name"
$ws1.Range("C1").Value = "04AM"
$ws1.Range("D1").Value = "05AM"
$ws1.Range("E1").Value = "06AM"
$ws1.Range("F1").Value = "07AM"
$ws1.Range("G1").Value = "08AM"
$ws1.Range("H1").Value = "09AM"
$ws1.Range("I1").Value = "10AM"
$ws1.Range("J1").Value = "11AM"
$ws1.Range("K1").Value = "Rates"

# Row 18 ("Total") needs the same A-column style (bold+border) as the other A-column cells.
# Copy the style from A17 before writing values so the copy does not clobber the new content.
$ws1.Range("A17").Copy($ws1.Range("A18"))

# Employee name labels (columns A and B) for rows 2-18
$ws1.Range("A2").Value = "Luis,Naula Jara"
$ws1.Range("B2").Value = "Luis,Naula Jara"
$ws1.Range("A3").Value = "Manuel,Martinez Suarez"
$ws1.Range("B3").Value = "Manuel,Martinez Suarez"
$ws1.Range("A4").Value = "Nancy,Ovillo"
$ws1.Range("B4").Value = "Nancy,Ovillo"
$ws1.Range("A5").Value = "Enmanuel,Vargas Rodriguez"
$ws1.Range("B5").Value = "Enmanuel,Vargas Rodriguez"
$ws1.Range("A6").Value = "Aurelio,Mercedes Hernandez"
$ws1.Range("B6").Value = "Aurelio,Mercedes Hernandez"
$ws1.Range("A7").Value = "Lorenzo,Peralta Santos"
$ws1.Range("B7").Value = "Lorenzo,Peralta Santos"
$ws1.Range("A8").Value = "Audrey,Henry"
$ws1.Range("B8").Value = "Audrey,Henry"
$ws1.Range("A9").Value = "Yave,Caba Corona"
$ws1.Range("B9").Value = "Yave,Caba Corona"
$ws1.Range("A10").Value = "Anthony,Alexander"
$ws1.Range("B10").Value = "Anthony,Alexander"
$ws1.Range("A11").Value = "Davis,Villavicencio Lenes"
$ws1.Range("B11").Value = "Davis,Villavicencio Lenes"
$ws1.Range("A12").Value = "Elia,Uruchima"
$ws1.Range("B12").Value = "Elia,Uruchima"
$ws1.Range("A13").Value = "Al,Smith"
$ws1.Range("B13").Value = "Al,Smith"
$ws1.Range("A14").Value = "Luis,Vargas Jaquez"
$ws1.Range("B14").Value = "Luis,Vargas Jaquez"
$ws1.Range("A15").Value = "Jonathan,Perez Henriquez"
$ws1.Range("B15").Value = "Jonathan,Perez Henriquez"
$ws1.Range("A16").Value = "Deron,Fuller"
$ws1.Range("B16").Value = "Deron,Fuller"
$ws1.Range("A17").Value = "Job-Bright,Dzameshie"
$ws1.Range("B17").Value = "Job-Bright,Dzameshie"
$ws1.Range("A18").Value = "Total"
$ws1.Range("B18").Value = "Total"

# Row 17 is now a regular employee row (no Rates/K value) - clear the stale Total-row carry-over
$ws1.Range("K17").ClearContents()

# Data grid (hours + rates) rows 2-18
$ws1.Range("C2").Value = 56
$ws1.Range("D2").Value = 168
$ws1.Range("E2").Value = 69
$ws1.Range("F2").Value = 9
$ws1.Range("G2").Value = 3
$ws1.Range("H2").Value = 19
$ws1.Range("I2").Value = 47
$ws1.Range("J2").Value = 0
$ws1.Range("K2").Value = 85
$ws1.Range("C3").Value = 75
$ws1.Range("D3").Value = 7
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 0
$ws1.Range("H3").Value = 0
$ws1.Range("I3").Value = 0
$ws1.Range("J3").Value = 0
$ws1.Range("K3").Value = 75
$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 1
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 40
$ws1.Range("G4").Value = 61
$ws1.Range("H4").Value = 0
$ws1.Range("I4").Value = 6
$ws1.Range("J4").Value = 118
$ws1.Range("K4").Value = 73
$ws1.Range("C5").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 109
$ws1.Range("F5").Value = 69
$ws1.Range("G5").Value = 64
$ws1.Range("H5").Value = 0
$ws1.Range("I5").Value = 0
$ws1.Range("J5").Value = 46
$ws1.Range("K5").Value = 72
$ws1.Range("C6").Value = 0
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 0
$ws1.Range("H6").Value = 0
$ws1.Range("I6").Value = 0
$ws1.Range("J6").Value = 69
$ws1.Range("K6").Value = 69
$ws1.Range("C7").Value = 28
$ws1.Range("D7").Value = 9
$ws1.Range("E7").Value = 6
$ws1.Range("F7").Value = 89
$ws1.Range("G7").Value = 74
$ws1.Range("H7").Value = 59
$ws1.Range("I7").Value = 0
$ws1.Range("J7").Value = 0
$ws1.Range("K7").Value = 62.5
$ws1.Range("C8").Value = 0
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("G8").Value = 27
$ws1.Range("H8").Value = 102
$ws1.Range("I8").Value = 89
$ws1.Range("J8").Value = 26
$ws1.Range("K8").Value = 61
$ws1.Range("C9").Value = 0
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 52
$ws1.Range("F9").Value = 0
$ws1.Range("G9").Value = 1
$ws1.Range("H9").Value = 72
$ws1.Range("I9").Value = 42
$ws1.Range("J9").Value = 0
$ws1.Range("K9").Value = 55.33333333333334
$ws1.Range("C10").Value = 43
$ws1.Range("D10").Value = 0
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 0
$ws1.Range("H10").Value = 0
$ws1.Range("I10").Value = 0
$ws1.Range("J10").Value = 0
$ws1.Range("K10").Value = 43
$ws1.Range("C11").Value = 35
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 0
$ws1.Range("H11").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("J11").Value = 0
$ws1.Range("K11").Value = 35
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 0
$ws1.Range("E12").Value = 0
$ws1.Range("F12").Value = 0
$ws1.Range("G12").Value = 0
$ws1.Range("H12").Value = 0
$ws1.Range("I12").Value = 0
$ws1.Range("J12").Value = 31
$ws1.Range("K12").Value = 31
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 0
$ws1.Range("E13").Value = 0
$ws1.Range("F13").Value = 2
$ws1.Range("G13").Value = 0
$ws1.Range("H13").Value = 0
$ws1.Range("I13").Value = 0
$ws1.Range("J13").Value = 0
$ws1.Range("C14").Value = 16
$ws1.Range("D14").Value = 0
$ws1.Range("E14").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("G14").Value = 0
$ws1.Range("H14").Value = 0
$ws1.Range("I14").Value = 0
$ws1.Range("J14").Value = 0
$ws1.Range("C15").Value = 1
$ws1.Range("D15").Value = 0
$ws1.Range("E15").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("G15").Value = 0
$ws1.Range("H15").Value = 0
$ws1.Range("I15").Value = 0
$ws1.Range("J15").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 0
$ws1.Range("E16").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("G16").Value = 0
$ws1.Range("H16").Value = 0
$ws1.Range("I16").Value = 0
$ws1.Range("J16").Value = 0
$ws1.Range("C17").Value = 19
$ws1.Range("D17").Value = 0
$ws1.Range("E17").Value = 0
$ws1.Range("F17").Value = 0
$ws1.Range("G17").Value = 0
$ws1.Range("H17").Value = 0
$ws1.Range("I17").Value = 0
$ws1.Range("J17").Value = 0
$ws1.Range("C18").Value = 278
$ws1.Range("D18").Value = 185
$ws1.Range("E18").Value = 236
$ws1.Range("F18").Value = 209
$ws1.Range("G18").Value = 230
$ws1.Range("H18").Value = 252
$ws1.Range("I18").Value = 184
$ws1.Range("J18").Value = 290
$ws1.Range("K18").Value = 226.5714285714286

# ---- Sheet "aux" ----
$ws2 = $wb.Worksheets.Item("aux")

# Header row (employee + total + hours label)
$ws2.Range("B1").Value = "Luis,Naula Jara"
$ws2.Range("C1").Value = "Manuel,Martinez Suarez"
$ws2.Range("D1").Value = "Nancy,Ovillo"
$ws2.Range("E1").Value = "Enmanuel,Vargas Rodriguez"
$ws2.Range("F1").Value = "Aurelio,Mercedes Hernandez"
$ws2.Range("G1").Value = "Lorenzo,Peralta Santos"
$ws2.Range("H1").Value = "Audrey,Henry"
$ws2.Range("I1").Value = "Yave,Caba Corona"
$ws2.Range("J1").Value = "Anthony,Alexander"
$ws2.Range("K1").Value = "Davis,Villavicencio Lenes"
$ws2.Range("L1").Value = "Elia,Uruchima"
$ws2.Range("M1").Value = "Al,Smith"
$ws2.Range("N1").Value = "Luis,Vargas Jaquez"
$ws2.Range("O1").Value = "Jonathan,Perez Henriquez"
$ws2.Range("P1").Value = "Deron,Fuller"
$ws2.Range("Q1").Value = "Job-Bright,Dzameshie"
$ws2.Range("R1").Value = "Total"
$ws2.Range("S1").Value = "hours"

# Hour labels (columns A and S) for rows 2-9
$ws2.Range("A2").Value = "04AM"
$ws2.Range("S2").Value = "04AM"
$ws2.Range("A3").Value = "05AM"
$ws2.Range("S3").Value = "05AM"
$ws2.Range("A4").Value = "06AM"
$ws2.Range("S4").Value = "06AM"
$ws2.Range("A5").Value = "07AM"
$ws2.Range("S5").Value = "07AM"
$ws2.Range("A6").Value = "08AM"
$ws2.Range("S6").Value = "08AM"
$ws2.Range("A7").Value = "09AM"
$ws2.Range("S7").Value = "09AM"
$ws2.Range("A8").Value = "10AM"
$ws2.Range("S8").Value = "10AM"
$ws2.Range("A9").Value = "11AM"
$ws2.Range("S9").Value = "11AM"

# Data grid rows 2-9
$ws2.Range("B2").Value = 56
$ws2.Range("C2").Value = 75
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 28
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 43
$ws2.Range("K2").Value = 35
$ws2.Range("L2").Value = 0
$ws2.Range("M2").Value = 0
$ws2.Range("N2").Value = 16
$ws2.Range("O2").Value = 1
$ws2.Range("P2").Value = 4
$ws2.Range("Q2").Value = 19
$ws2.Range("R2").Value = 278
$ws2.Range("B3").Value = 168
$ws2.Range("C3").Value = 7
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 9
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0
$ws2.Range("L3").Value = 0
$ws2.Range("M3").Value = 0
$ws2.Range("N3").Value = 0
$ws2.Range("O3").Value = 0
$ws2.Range("P3").Value = 0
$ws2.Range("Q3").Value = 0
$ws2.Range("R3").Value = 185
$ws2.Range("B4").Value = 69
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 109
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 6
$ws2.Range("H4").Value = 0
$ws2.Range("I4").Value = 52
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0
$ws2.Range("L4").Value = 0
$ws2.Range("M4").Value = 0
$ws2.Range("N4").Value = 0
$ws2.Range("O4").Value = 0
$ws2.Range("P4").Value = 0
$ws2.Range("Q4").Value = 0
$ws2.Range("R4").Value = 236
$ws2.Range("B5").Value = 9
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 40
$ws2.Range("E5").Value = 69
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 89
$ws2.Range("H5").Value = 0
$ws2.Range("I5").Value = 0
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0
$ws2.Range("L5").Value = 0
$ws2.Range("M5").Value = 2
$ws2.Range("N5").Value = 0
$ws2.Range("O5").Value = 0
$ws2.Range("P5").Value = 0
$ws2.Range("Q5").Value = 0
$ws2.Range("R5").Value = 209
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 61
$ws2.Range("E6").Value = 64
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 74
$ws2.Range("H6").Value = 27
$ws2.Range("I6").Value = 1
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0
$ws2.Range("L6").Value = 0
$ws2.Range("M6").Value = 0
$ws2.Range("N6").Value = 0
$ws2.Range("O6").Value = 0
$ws2.Range("P6").Value = 0
$ws2.Range("Q6").Value = 0
$ws2.Range("R6").Value = 230
$ws2.Range("B7").Value = 19
$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 0
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 59
$ws2.Range("H7").Value = 102
$ws2.Range("I7").Value = 72
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0
$ws2.Range("L7").Value = 0
$ws2.Range("M7").Value = 0
$ws2.Range("N7").Value = 0
$ws2.Range("O7").Value = 0
$ws2.Range("P7").Value = 0
$ws2.Range("Q7").Value = 0
$ws2.Range("R7").Value = 252
$ws2.Range("B8").Value = 47
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 6
$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 0
$ws2.Range("H8").Value = 89
$ws2.Range("I8").Value = 42
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0
$ws2.Range("L8").Value = 0
$ws2.Range("M8").Value = 0
$ws2.Range("N8").Value = 0
$ws2.Range("O8").Value = 0
$ws2.Range("P8").Value = 0
$ws2.Range("Q8").Value = 0
$ws2.Range("R8").Value = 184
$ws2.Range("B9").Value = 0
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 118
$ws2.Range("E9").Value = 46
$ws2.Range("F9").Value = 69
$ws2.Range("G9").Value = 0
$ws2.Range("H9").Value = 26
$ws2.Range("I9").Value = 0
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0
$ws2.Range("L9").Value = 31
$ws2.Range("M9").Value = 0
$ws2.Range("N9").Value = 0
$ws2.Range("O9").Value = 0
$ws2.Range("P9").Value = 0
$ws2.Range("Q9").Value = 0
$ws2.Range("R9").Value = 290

# ---- Sheet "floors_table" ----
$ws3 = $wb.Worksheets.Item("floors_table")

# Header row
$ws3.Range("B1").Value = "level"
$ws3.Range("C1").Value = "Night shift"
$ws3.Range("D1").Value = "Morning shift"
$ws3.Range("E1").Value = "Afternoon shift"
$ws3.Range("F1").Value = "Total"

# Level labels column B, rows 2-5
$ws3.Range("B2").Value = "A1"
$ws3.Range("B3").Value = "A2"
$ws3.Range("B4").Value = "A3"
$ws3.Range("B5").Value = "A4"

# Data grid rows 2-5
$ws3.Range("C2").Value = 246
$ws3.Range("D2").Value = 272
$ws3.Range("E2").Value = 194
$ws3.Range("F2").Value = 712
$ws3.Range("C3").Value = 301
$ws3.Range("D3").Value = 382
$ws3.Range("E3").Value = 376
$ws3.Range("F3").Value = 1059
$ws3.Range("C4").Value = 643
$ws3.Range("D4").Value = 661
$ws3.Range("E4").Value = 683
$ws3.Range("F4").Value = 1987
$ws3.Range("C5").Value = 498
$ws3.Range("D5").Value = 549
$ws3.Range("E5").Value = 456
$ws3.Range("F5").Value = 1503

# ---- Sheet "results_table" ----
$ws4 = $wb.Worksheets.Item("results_table")

# Header row
$ws4.Range("B1").Value = "Expected Results"
$ws4.Range("C1").Value = "Net Results"
$ws4.Range("D1").Value = "Difference"

# Data row 2
$ws4.Range("B2").Value = 3500
$ws4.Range("C2").Value = 1864
$ws4.Range("D2").Value = -1636

